# Insert a new "date" question row into the XLSForm "survey" sheet right
# after the "Monitoring visits" select_one question (row 8) and before
# "Name of monitoring visitor" (previously row 9), shifting the remaining
# rows of the group down by one and adding a new shared-string entries for
# date / _2_201 / Date of the monitoring visit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Shift rows 9..12 down to 10..13, opening up a blank row 9.
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row with the new "date" field.
$ws.Cells.Item(9, 1).Value = "date"
$ws.Cells.Item(9, 2).Value = "_2_201"
$ws.Cells.Item(9, 3).Value = "Date of the monitoring visit"

# Match the saved selection state from the authored edit.
[void]$ws.Range("C9").Select()
